# Applies a cyclic rotation of the weekly price blocks:
#   rows 6-7 (date 44223) -> rows 2-3
#   rows 2-3 (date 44559) -> rows 4-5
#   rows 4-5 (date 44574) -> rows 6-7
# Only the data values change; cell styles/formatting are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values after rotation, keyed by row number.
$targets = @{
    2 = @{ D = 44223; M = 100; N = 3500; O = 4000; P = 3750; S = 1875 }
    3 = @{ D = 44223; M = 50;  N = 3000; O = 3000; P = 3000; S = 1500 }
    4 = @{ D = 44559; M = 200; N = 6000; O = 7000; P = 6500; S = 3250 }
    5 = @{ D = 44559; M = 100; N = 5000; O = 5000; P = 5000; S = 2500 }
    6 = @{ D = 44574; M = 200; N = 6000; O = 7000; P = 6500; S = 3250 }
    7 = @{ D = 44574; M = 100; N = 5000; O = 5000; P = 5000; S = 2500 }
}

foreach ($row in $targets.Keys) {
    $vals = $targets[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}
